$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.714.27'
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").Value = '3.636.41'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '2.71'
$ws.Range("E4").Value = '  +19.64%  '
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '224.96'
$ws.Range("E6").Value = '  -5.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '645.84'
$ws.Range("E7").Value = '  -1.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.424'
$ws.Range("E8").Value = '  -4.31%  '
$ws.Range("E9").Value = '  +3.50%  '
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").Value = '3.632.05'
$ws.Range("E11").Value = '  -2.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '51.68'
$ws.Range("E12").Value = '  +13.63%  '
$ws.Range("E13").Value = '  +5.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000297'
$ws.Range("E14").Value = '  -4.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.56'
$ws.Range("E15").Value = '  -4.23%  '
$ws.Range("D16").Value = '4.316.26'
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.85'
$ws.Range("E17").Value = '  +30.78%  '
$ws.Range("D18").Value = '95.492.69'
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.86'
$ws.Range("E19").Value = '  -5.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.88'
$ws.Range("E20").Value = '  +6.05%  '
$ws.Range("D21").Value = '3.632.36'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.315'
$ws.Range("E22").Value = '  +48.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.538'
$ws.Range("E23").Value = '  -2.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.36'
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.20'
$ws.Range("E27").Value = '  +4.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000204'
$ws.Range("E28").Value = '  -8.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.70'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").Value = '3.803.38'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.74'
$ws.Range("E31").Value = '  +6.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("E32").Value = '  +5.78%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '33.86'
$ws.Range("E36").Value = '  +2.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.184'
$ws.Range("E37").Value = '  -4.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0566'
$ws.Range("E39").Value = '  +22.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.62'
$ws.Range("E40").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.27'
$ws.Range("E43").Value = '  +4.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.508'
$ws.Range("E44").Value = '  +3.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.03'
$ws.Range("E45").Value = '  +5.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.57'
$ws.Range("E46").Value = '  +1.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.04'
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.158'
$ws.Range("E48").Value = '  -6.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.39'
$ws.Range("E49").Value = '  +4.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '237.65'
$ws.Range("E50").Value = '  +13.17%  '
$ws.Range("E51").Value = '  -0.89%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '536.36'
$ws.Range("E24").Value = '  +0.97%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.10'
$ws.Range("E25").Value = '  +12.82%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  +4.08%  '

$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.646'
$ws.Range("E35").Value = '  +6.81%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '604.33'
$ws.Range("E41").Value = '  -5.64%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.02%  '

